# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml -> currently "Office Theme" (used only by the Notes Master)
#   ppt/theme/theme2.xml -> currently "Integral" / "Red Violet" (used by the Slide Master
#                             and the presentation itself)
#
# The authored edit swaps the *content* of those two theme parts, so the
# Slide Master ends up on the plain "Office Theme" palette and the Notes
# Master ends up on the old "Integral" / "Red Violet" palette.
#
# The only part of that swap that is reachable through the PowerPoint object
# model is the colour scheme that is actually applied to the deck (i.e. the
# theme driving the Slide Master / Presentation, which lives in theme2.xml).
# We recolor it to the stock Office theme palette, which is exactly the
# palette the target theme2.xml ends up with.

function BgrVal($r, $g, $b) { return ($b * 65536) + ($g * 256) + $r }

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

# Office theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# which is what theme2.xml ends up containing after the swap.
$cs.Colors(1).RGB  = BgrVal 0x00 0x00 0x00   # dk1      000000
$cs.Colors(2).RGB  = BgrVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$cs.Colors(3).RGB  = BgrVal 0x44 0x54 0x6A   # dk2      44546A
$cs.Colors(4).RGB  = BgrVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$cs.Colors(5).RGB  = BgrVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$cs.Colors(6).RGB  = BgrVal 0xED 0x7D 0x31   # accent2  ED7D31
$cs.Colors(7).RGB  = BgrVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$cs.Colors(8).RGB  = BgrVal 0xFF 0xC0 0x00   # accent4  FFC000
$cs.Colors(9).RGB  = BgrVal 0x44 0x72 0xC4   # accent5  4472C4
$cs.Colors(10).RGB = BgrVal 0x70 0xAD 0x47   # accent6  70AD47
$cs.Colors(11).RGB = BgrVal 0x05 0x63 0xC1   # hlink    0563C1
$cs.Colors(12).RGB = BgrVal 0x95 0x4F 0x72   # folHlink 954F72
